# Fix list level numbering: decrement the indent level of list paragraphs
# by one, since PowerPoint's top-level list content should sit at the same
# level as top-level paragraph content (only the list style differs).
#
# This mirrors the pandoc pptx writer fix: only continuation paragraphs of
# a list should have their level incremented, not every list paragraph.

$p = $ppt.ActivePresentation

# --- Slide 1: "Bulleted bulleted lists." example -----------------------
$s1 = $p.Slides.Item(1)
$body1 = $s1.Shapes.Item(2).TextFrame.TextRange

for ($i = 1; $i -le $body1.Paragraphs().Count; $i++) {
    $para = $body1.Paragraphs($i, 1)
    $para.IndentLevel = $para.IndentLevel - 1
}

# --- Slide 2: "Lists can also be numbered" example ----------------------
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2).TextFrame.TextRange

# The first paragraph on this slide ("Lists can also be numbered:") is not
# part of the list and keeps its level; only the numbered-list paragraphs
# (2 through 6) shift down a level.
for ($i = 2; $i -le $body2.Paragraphs().Count; $i++) {
    $para = $body2.Paragraphs($i, 1)
    $para.IndentLevel = $para.IndentLevel - 1
}
